# daily auto push: 2026-01-28 18:52 UTC
# Two new data rows (2026/01/28 water-day and 2026/01/29 thursday) were
# inserted into the daily log right after the existing 2026/01/28 block
# (old row 739), pushing every subsequent row (old 740..781) down by two
# rows (new 742..783).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows at 740-741; everything at/below shifts down two rows
# (matches Excel's default whole-row Insert behaviour, i.e. Shift:=xlShiftDown).
$ws.Rows("740:741").Insert()

# Column A holds the date as plain text (e.g. "2026/01/28"), not a real
# date serial, matching every other row in the sheet - so force text
# formatting before writing, then restore the default style so the new
# cells match the rest of the column.
$ws.Range("A740:A741").NumberFormat = "@"

$ws.Range("A740").Value = "2026/01/28"
$ws.Range("B740").Value = "水"
$ws.Range("C740").Value = 23
$ws.Range("D740").Value = 201

$ws.Range("A741").Value = "2026/01/29"
$ws.Range("B741").Value = "木"
$ws.Range("C741").Value = 2
$ws.Range("D741").Value = 201

$ws.Range("A740:A741").Style = "Normal"
